$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Remove curly quotes around $param in the two ACTION template strings.
$ws.Range("G7").Value = "`$p.setModifiedRatingFactor(`$param);"
$ws.Range("H7").Value = "`$p.setTier(`$param);"

# Wrap the tier literal values in curly quotes.
$ws.Range("H9").Value = [char]0x201C + "Tier1" + [char]0x201D
$ws.Range("H10").Value = [char]0x201C + "Tier1" + [char]0x201D
$ws.Range("H11").Value = [char]0x201C + "Tier2" + [char]0x201D
$ws.Range("H12").Value = [char]0x201C + "Tier2" + [char]0x201D
$ws.Range("H13").Value = [char]0x201C + "Tier3" + [char]0x201D
$ws.Range("H14").Value = [char]0x201C + "Tier3" + [char]0x201D
$ws.Range("H15").Value = [char]0x201C + "Tier4" + [char]0x201D
$ws.Range("H16").Value = [char]0x201C + "Tier4" + [char]0x201D
$ws.Range("H17").Value = [char]0x201C + "Tier5" + [char]0x201D
$ws.Range("H18").Value = [char]0x201C + "Tier5" + [char]0x201D

# Update the active selection on the sheet (was F29, now I9).
$ws.Range("I9").Select()
